# Refresh crypto price/volume data to match the latest scrape.
# A handful of rows (33/34 and 36/37) also changed which coin occupies
# that row, so Coin/Link/Price are rewritten for those too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.923.10'
$ws.Range('E2').Value = '  +3.09%  '
# Row 3
$ws.Range('D3').Value = '3.441.57'
$ws.Range('E3').Value = '  +1.79%  '
# Row 4
$ws.Range('E4').Value = '  -0.13%  '
# Row 5
$ws.Range('D5').Value = "'570.44"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.17%  '
# Row 6
$ws.Range('D6').Value = "'184.79"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.56%  '
# Row 7
$ws.Range('D7').Value = "'0.634"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.66%  '
# Row 8
$ws.Range('D8').Value = '3.436.41'
$ws.Range('E8').Value = '  +1.93%  '
# Row 9
$ws.Range('D9').Value = "'1.00"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '
# Row 10
$ws.Range('D10').Value = "'0.178"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.96%  '
# Row 11
$ws.Range('D11').Value = "'0.644"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.94%  '
# Row 12
$ws.Range('D12').Value = "'55.37"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.34%  '
# Row 13
$ws.Range('E13').Value = '  +1.92%  '
# Row 14
$ws.Range('D14').Value = "'9.38"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.09%  '
# Row 15
$ws.Range('D15').Value = '3.986.56'
$ws.Range('E15').Value = '  +1.35%  '
# Row 16
$ws.Range('D16').Value = "'18.56"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.55%  '
# Row 17
$ws.Range('D17').Value = '3.445.95'
$ws.Range('E17').Value = '  +1.69%  '
# Row 18
$ws.Range('E18').Value = '  +0.51%  '
# Row 19
$ws.Range('D19').Value = '66.842.82'
$ws.Range('E19').Value = '  +2.87%  '
# Row 20
$ws.Range('D20').Value = "'12.01"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.86%  '
# Row 21
$ws.Range('E21').Value = '  +1.89%  '
# Row 22
$ws.Range('D22').Value = "'475.32"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.64%  '
# Row 23
$ws.Range('D23').Value = "'4.98"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.06%  '
# Row 24
$ws.Range('D24').Value = "'14.96"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +10.74%  '
# Row 25
$ws.Range('D25').Value = "'4.19"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.49%  '
# Row 26
$ws.Range('D26').Value = "'89.61"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.68%  '
# Row 27
$ws.Range('D27').Value = "'2.96"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.12%  '
# Row 28
$ws.Range('D28').Value = "'10.97"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.16%  '
# Row 29
$ws.Range('D29').Value = "'8.92"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.22%  '
# Row 30
$ws.Range('D30').Value = "'31.53"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.38%  '
# Row 31
$ws.Range('E31').Value = '  +3.27%  '
# Row 32
$ws.Range('D32').Value = "'11.63"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.31%  '
# Row 33
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = "'588.76"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.24%  '
# Row 34
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = "'63.08"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.91%  '
# Row 35
$ws.Range('E35').Value = '  +1.63%  '
# Row 36
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = "'0.999"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.09%  '
# Row 37
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = "'0.148"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.83%  '
# Row 38
$ws.Range('E38').Value = '  +0.50%  '
# Row 39
$ws.Range('D39').Value = "'0.392"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.13%  '
# Row 40
$ws.Range('D40').Value = "'36.60"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.00%  '
# Row 41
$ws.Range('D41').Value = '0.0₃0774'
$ws.Range('E41').Value = '  +4.21%  '
# Row 42
$ws.Range('D42').Value = '3.129.83'
$ws.Range('E42').Value = '  +1.56%  '
# Row 43
$ws.Range('D43').Value = "'2.91"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.83%  '
# Row 44
$ws.Range('D44').Value = "'2.63"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.56%  '
# Row 45
$ws.Range('D45').Value = "'0.0425"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.28%  '
# Row 46
$ws.Range('E46').Value = '  +19.94%  '
# Row 47
$ws.Range('D47').Value = "'3.22"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.45%  '
# Row 48
$ws.Range('E48').Value = '  +0.36%  '
# Row 49
$ws.Range('D49').Value = "'0.999"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.17%  '
# Row 50
$ws.Range('D50').Value = "'142.10"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.96%  '
# Row 51
$ws.Range('D51').Value = "'8.69"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.63%  '
